$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.582.05'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '1.877.13'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.95'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2834'
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06532'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').Value = '1.917.90'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07481'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.096'
$ws.Range('E13').Value = '  -1.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.69'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6635'
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('D16').Value = '30.543.12'
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.33'
$ws.Range('E17').Value = '  -2.06%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007615'
$ws.Range('E19').Value = '  -1.42%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.89'
$ws.Range('E20').Value = '  +17.49%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.205.60'
$ws.Range('E21').Value = '  +3.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.310'
$ws.Range('E22').Value = '  -2.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.204'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.310'
$ws.Range('E25').Value = '  -1.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.35'
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.71'
$ws.Range('E27').Value = '  +2.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.952'
$ws.Range('E28').Value = '  +1.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.461'
$ws.Range('E29').Value = '  +1.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09534'
$ws.Range('E30').Value = '  +4.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.311'
$ws.Range('E31').Value = '  +0.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.032'
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05026'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.215'
$ws.Range('E34').Value = '  +6.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7481'
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01837'
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.075'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9108'
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '106.08'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.813'
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4275'
$ws.Range('E43').Value = '  -1.09%  '
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.469'
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.52'
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('E47').Value = '  -4.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.477'
$ws.Range('E48').Value = '  -6.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.918'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.84'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3883'
$ws.Range('E51').Value = '  +0.35%  '
